# Apply the scraped-price refresh for the cryptos worksheet.
# For each changed cell we set the literal text; numeric-looking values
# (e.g. "1.000", "50.10") are pushed through a Text number format so
# Excel does not silently re-interpret them as numbers, then the
# temporary formatting is cleared so the cell style is left untouched
# (matching the source workbook, where these are plain General cells).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "23.444.83"
$ws.Range("E2").Value = "  -1.29%  "

# Row 3
$ws.Range("D3").Value = "1.645.92"
$ws.Range("E3").Value = "  -0.59%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.21%  "

# Row 5
$ws.Range("E5").Value = "  -0.04%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "299.19"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -1.49%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3789"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.24%  "

# Row 8
$ws.Range("B8").Value = "OKB"
$ws.Range("C8").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "50.10"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -1.93%  "

# Row 9
$ws.Range("B9").Value = "Cardano"
$ws.Range("C9").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3515"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -3.17%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08082"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -1.76%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.214"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -3.29%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.001"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -0.19%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.12"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -2.53%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.374"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -2.60%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.327"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -2.03%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001202"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -3.29%  "

# Row 17
$ws.Range("D17").Value = "1.648.40"
$ws.Range("E17").Value = "  -0.57%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "96.51"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -1.34%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06999"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.22%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.720"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -1.27%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.37"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -2.33%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.000"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.09%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.36"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -3.71%  "

# Row 24
$ws.Range("D24").Value = "23.467.89"
$ws.Range("E24").Value = "  -1.18%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.484"
$ws.Range("D25").ClearFormats()

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.917"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -5.20%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.85"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -2.29%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "153.25"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +1.07%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.198"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.48%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "132.72"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -1.34%  "

# Row 31
$ws.Range("D31").Value = "1.827.53"
$ws.Range("E31").Value = "  -0.70%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.887"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.45%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.122"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -3.07%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.41"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -3.52%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9840"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -8.84%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02706"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -4.35%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.08747"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.87%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2432"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -3.79%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.920"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -3.46%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06812"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -4.19%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "12.90"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -2.79%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6878"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -2.80%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.297"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -3.43%  "

# Row 44
$ws.Range("E44").Value = "  -2.58%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.000"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.04%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6353"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -3.13%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.254"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -3.35%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.906"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -1.40%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.07713"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -3.17%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "127.17"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.82%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.145"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -4.16%  "

